# Updates Leve profitability figures across the Goblin_Profits sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect refreshed market
# board pricing pulled in by the scheduled runner.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# row 5 (Leve Item ID 5503)
$ws.Range("H5").Value = 161.28572
$ws.Range("I5").Value = 112.25
$ws.Range("K5").Value = 112.25
$ws.Range("M5").Value = 2.75

# row 9 (Leve Item ID 5487)
$ws.Range("H9").Value = 82.35714
$ws.Range("I9").Value = 62.4
$ws.Range("K9").Value = 62.4
$ws.Range("M9").Value = 106.6

# row 12 (Leve Item ID 5515)
$ws.Range("H12").Value = 181.33333
$ws.Range("I12").Value = 139.75
$ws.Range("J12").Value = 264.5
$ws.Range("K12").Value = 139.75
$ws.Range("L12").Value = 264.5
$ws.Range("M12").Value = 30.25
$ws.Range("N12").Value = -604.5

# row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 3876.7778
$ws.Range("J40").Value = 4127.7144
$ws.Range("L40").Value = 4127.7144
$ws.Range("N40").Value = -4477.7144

# row 64 (Leve Item ID 5506)
$ws.Range("H64").Value = 8249.385
$ws.Range("I64").Value = 5748.5
$ws.Range("J64").Value = 9812.4375
$ws.Range("K64").Value = 5748.5
$ws.Range("L64").Value = 9812.4375
$ws.Range("M64").Value = -5500.5
$ws.Range("N64").Value = -10308.4375

# row 67 (Leve Item ID 5506)
$ws.Range("H67").Value = 8249.385
$ws.Range("I67").Value = 5748.5
$ws.Range("J67").Value = 9812.4375
$ws.Range("K67").Value = 5748.5
$ws.Range("L67").Value = 9812.4375
$ws.Range("M67").Value = -4890.5
$ws.Range("N67").Value = -11528.4375

# row 106 (Leve Item ID 19903)
$ws.Range("H106").Value = 3090.7368
$ws.Range("I106").Value = 2836.7646
$ws.Range("K106").Value = 2836.7646
$ws.Range("M106").Value = -2205.7646

# row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 1971.4
$ws.Range("J137").Value = 1500
$ws.Range("L137").Value = 4500
$ws.Range("N137").Value = -9600


# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# row 24 (Leve Item ID 18363)
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

# row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 3919.6667
$ws.Range("I32").Value = 3866.5898
$ws.Range("K32").Value = 3866.5898
$ws.Range("M32").Value = -3579.5898

# row 94 (Leve Item ID 18055)
$ws.Range("H94").Value = 24249.25
$ws.Range("J94").Value = 24249.25
$ws.Range("L94").Value = 24249.25
$ws.Range("N94").Value = -26051.25

# row 95 (Leve Item ID 18204)
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492

# row 98 (Leve Item ID 18371)
$ws.Range("H98").Value = 21675
$ws.Range("J98").Value = 21675
$ws.Range("L98").Value = 21675
$ws.Range("N98").Value = -27665

# row 100 (Leve Item ID 18363)
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

# row 112 (Leve Item ID 25808)
$ws.Range("H112").Value = 38000
$ws.Range("J112").Value = 38000
$ws.Range("L112").Value = 38000
$ws.Range("N112").Value = -40954


# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# row 82 (Leve Item ID 11877)
$ws.Range("H82").Value = 18436.428
$ws.Range("I82").Value = 11509.167
$ws.Range("J82").Value = 60000
$ws.Range("K82").Value = 11509.167
$ws.Range("L82").Value = 60000
$ws.Range("M82").Value = -11126.167
$ws.Range("N82").Value = -60766

# row 85 (Leve Item ID 11877)
$ws.Range("H85").Value = 18436.428
$ws.Range("I85").Value = 11509.167
$ws.Range("J85").Value = 60000
$ws.Range("K85").Value = 11509.167
$ws.Range("L85").Value = 60000
$ws.Range("M85").Value = -10183.167
$ws.Range("N85").Value = -62652


# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# row 2 (Leve Item ID 1820)
$ws.Range("H2").Value = 341.85715
$ws.Range("I2").Value = 358.6
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 358.6
$ws.Range("L2").Value = 300
$ws.Range("M2").Value = -245.6
$ws.Range("N2").Value = -526

# row 22 (Leve Item ID 5367)
$ws.Range("H22").Value = 1673.9231
$ws.Range("I22").Value = 946
$ws.Range("J22").Value = 2297.8572
$ws.Range("K22").Value = 946
$ws.Range("L22").Value = 2297.8572
$ws.Range("M22").Value = -596
$ws.Range("N22").Value = -2997.8572

# row 43 (Leve Item ID 18504)
$ws.Range("H43").Value = 20737.5
$ws.Range("J43").Value = 20737.5
$ws.Range("L43").Value = 20737.5
$ws.Range("N43").Value = -21105.5

# row 101 (Leve Item ID 18504)
$ws.Range("H101").Value = 20737.5
$ws.Range("J101").Value = 20737.5
$ws.Range("L101").Value = 20737.5
$ws.Range("N101").Value = -27227.5


# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# row 16 (Leve Item ID 4641)
$ws.Range("H16").Value = 330.875
$ws.Range("J16").Value = 337
$ws.Range("L16").Value = 1011
$ws.Range("N16").Value = -1357

# row 93 (Leve Item ID 19808)
$ws.Range("H93").Value = 11474.7
$ws.Range("I93").Value = 993
$ws.Range("J93").Value = 18462.5
$ws.Range("K93").Value = 2979
$ws.Range("L93").Value = 55387.5
$ws.Range("M93").Value = -1107
$ws.Range("N93").Value = -59131.5

# row 97 (Leve Item ID 19846)
$ws.Range("H97").Value = 306.0625
$ws.Range("J97").Value = 338.5
$ws.Range("L97").Value = 1015.5
$ws.Range("N97").Value = -2007.5

# row 100 (Leve Item ID 19831)
$ws.Range("H100").Value = 9159.5
$ws.Range("I100").Value = 596
$ws.Range("K100").Value = 1788
$ws.Range("M100").Value = -977

# row 107 (Leve Item ID 27838)
$ws.Range("H107").Value = 3290.9092
$ws.Range("I107").Value = 4941.8335
$ws.Range("K107").Value = 14825.5005
$ws.Range("M107").Value = -12905.5005


# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# row 7 (Leve Item ID 4197)
$ws.Range("H7").Value = 1893500
$ws.Range("J7").Value = 5583.3335
$ws.Range("L7").Value = 5583.3335
$ws.Range("N7").Value = -5807.3335

# row 8 (Leve Item ID 4197)
$ws.Range("H8").Value = 1893500
$ws.Range("J8").Value = 5583.3335
$ws.Range("L8").Value = 5583.3335
$ws.Range("N8").Value = -5861.3335

# row 11 (Leve Item ID 4422)
$ws.Range("H11").Value = 13466167
$ws.Range("J11").Value = 22000
$ws.Range("L11").Value = 22000
$ws.Range("N11").Value = -22278

# row 39 (Leve Item ID 18264)
$ws.Range("H39").Value = 50001
$ws.Range("J39").Value = 50001
$ws.Range("L39").Value = 50001
$ws.Range("N39").Value = -51065

# row 62 (Leve Item ID 11983)
$ws.Range("H62").Value = 49250
$ws.Range("J62").Value = 49250
$ws.Range("L62").Value = 49250
$ws.Range("N62").Value = -50622

# row 65 (Leve Item ID 11983)
$ws.Range("H65").Value = 49250
$ws.Range("J65").Value = 49250
$ws.Range("L65").Value = 147750
$ws.Range("N65").Value = -154614

# row 101 (Leve Item ID 18513)
$ws.Range("H101").Value = 18750
$ws.Range("J101").Value = 18750
$ws.Range("L101").Value = 18750
$ws.Range("N101").Value = -25240


# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# row 3 (Leve Item ID 3537)
$ws.Range("H3").Value = 4604
$ws.Range("J3").Value = 4604
$ws.Range("L3").Value = 4604
$ws.Range("N3").Value = -4828

# row 4 (Leve Item ID 3788)
$ws.Range("H4").Value = 5672.6665
$ws.Range("I4").Value = 5509
$ws.Range("K4").Value = 5509
$ws.Range("M4").Value = -5396

# row 15 (Leve Item ID 3537)
$ws.Range("H15").Value = 4604
$ws.Range("J15").Value = 4604
$ws.Range("L15").Value = 4604
$ws.Range("N15").Value = -4944

# row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 3558.8823
$ws.Range("I22").Value = 3708.1667
$ws.Range("J22").Value = 3477.4546
$ws.Range("K22").Value = 3708.1667
$ws.Range("L22").Value = 3477.4546
$ws.Range("M22").Value = -3413.1667
$ws.Range("N22").Value = -4067.4546

# row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 3558.8823
$ws.Range("I27").Value = 3708.1667
$ws.Range("J27").Value = 3477.4546
$ws.Range("K27").Value = 3708.1667
$ws.Range("L27").Value = 3477.4546
$ws.Range("M27").Value = -3601.1667
$ws.Range("N27").Value = -3691.4546

# row 28 (Leve Item ID 3788)
$ws.Range("H28").Value = 5672.6665
$ws.Range("I28").Value = 5509
$ws.Range("K28").Value = 5509
$ws.Range("M28").Value = -5277

# row 37 (Leve Item ID 3788)
$ws.Range("H37").Value = 5672.6665
$ws.Range("I37").Value = 5509
$ws.Range("K37").Value = 5509
$ws.Range("M37").Value = -5402

# row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 2781
$ws.Range("J46").Value = 3681.7
$ws.Range("L46").Value = 3681.7
$ws.Range("N46").Value = -4057.7


# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# row 9 (Leve Item ID 3015)
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()

# row 92 (Leve Item ID 18088)
$ws.Range("H92").Value = 40000
$ws.Range("J92").Value = 40000
$ws.Range("L92").Value = 40000
$ws.Range("N92").Value = -44992

# row 95 (Leve Item ID 18243)
$ws.Range("H95").Value = 32567.8
$ws.Range("J95").Value = 32567.8
$ws.Range("L95").Value = 32567.8
$ws.Range("N95").Value = -38059.8

